# Apply the Maximum Capacity Factor workbook update:
#  - Update the "About" sheet's last-updated date (C1)
#  - Update the MCF sheet's capacity-factor values from 0.85/0.95 to 1
#  - Leave the active-cell selection on MCF at B17 (matches saved view state)

$wb = $excel.ActiveWorkbook

# --- About sheet: bump the last updated date ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45392

# --- MCF sheet: set capacity factors to 1 for the affected plant types ---
$wsMCF = $wb.Worksheets.Item("MCF")

$cells = @("B2", "B3", "B4", "B6", "B10", "B11", "B12", "B13", "B14", "B16", "B17", "B18")
foreach ($cell in $cells) {
    $wsMCF.Range($cell).Value = 1
}

# Formula-driven cells (B19, B20, B21, B22, B24, B25) recalc automatically
# since they reference B2, B4, B10, B14.
$wb.Application.Calculate()

# Restore the saved selection/active cell on the MCF sheet.
$wsMCF.Activate()
$wsMCF.Range("B17").Select()
